$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly-added day's data on row 25
$ws.Range("A25").Value = 45971
$ws.Range("B25").Value = 598
$ws.Range("C25").Value = 20
$ws.Range("D25").Value = 578

# Update the sheet view: selection moves to B28 and the view scrolls so
# row 17 is the top-left visible row (topLeftCell="A17")
$ws.Range("B28").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 17
$win.ScrollColumn = 1
